$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ---
# Row 2: fuel changed from GAS_NATURAL to DIESEL, date moved from Apr-22 to Oct-22
$ws.Range("B2").Value = "DIESEL"
$ws.Range("E2").Value = 44856
$ws.Range("E2").NumberFormat = "mmm-yy"

# Rows 3-6: date moved from May-22 to Nov-22
$ws.Range("E3").Value = 44887
$ws.Range("E3").NumberFormat = "mmm-yy"
$ws.Range("E4").Value = 44887
$ws.Range("E4").NumberFormat = "mmm-yy"
$ws.Range("E5").Value = 44887
$ws.Range("E5").NumberFormat = "mmm-yy"
$ws.Range("E6").Value = 44887
$ws.Range("E6").NumberFormat = "mmm-yy"

# --- Append new COMBUSTION_FIJA rows for other fuel types ---
$fuels = @("KEROSENE", "FUEL_OIL", "NAFTA", "CARBON_DE_LEÑA", "LEÑA")
$row = 7
foreach ($fuel in $fuels) {
    $ws.Cells.Item($row, 1).Value = "COMBUSTION_FIJA"
    $ws.Cells.Item($row, 2).Value = $fuel
    $ws.Cells.Item($row, 3).Value = 50
    $ws.Cells.Item($row, 4).Value = "MENSUAL"
    $ws.Cells.Item($row, 5).Value = 44856
    $ws.Cells.Item($row, 5).NumberFormat = "mmm-yy"
    $row++
}

# --- Autofit columns A:E to match the resulting best-fit widths ---
$ws.Columns("A:E").AutoFit()

# --- Update selection to match final state ---
$ws.Range("G12").Select()
